$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.307.99"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "1.883.50"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'236.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'0.4836"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").Value = "'0.2883"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "'0.06597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "1.875.74"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").Value = "'16.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "'0.07332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'5.124"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'87.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.6593"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "30.261.99"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9994"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007750"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "'5.406"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "2.125.12"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "'0.9990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "'195.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.40%  "
$ws.Range("D24").Value = "'6.166"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'9.273"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.51%  "
$ws.Range("D26").Value = "'164.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "'18.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").Value = "'1.434"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'4.303"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'0.09148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'4.025"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'0.05064"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.134"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7213"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.80%  "
$ws.Range("D36").Value = "'2.695"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").Value = "'0.01783"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("D38").Value = "'2.633"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "'0.9186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "'2.053"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "'105.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'0.4299"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.23%  "
$ws.Range("D43").Value = "'5.807"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "'7.468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").Value = "'0.1314"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.39%  "
$ws.Range("D47").Value = "'64.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.82%  "
$ws.Range("D48").Value = "'1.528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.68%  "
$ws.Range("D49").Value = "'8.906"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "'33.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("D51").Value = "'0.05738"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.76%  "
